# =====================================================================
# Add "2022-Q4" quarterly holdings sheet and update the summary ("总计")
# sheet to account for it, matching the upstream commit "feat: add
# 2022-Q4 data".
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet ("总计"): insert a new row right under the header for
#    the 2022-Q4 totals, shifting every existing quarter down by one row
#    and keeping the running index (column A) sequential.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# New row 2 holds the 2022-Q4 totals.
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 29
$summary.Cells.Item(2, 4).Value = 23.18

# The row Insert() leaves B2:D2 carrying stray header-row formatting and
# A2 with no formatting at all; line them up with their neighbours.
$summary.Range("B2:D2").ClearFormats()
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)

# Renumber the running index in column A for the rows that shifted down.
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) New "2022-Q4" sheet: clone the "2022-Q3" sheet (so fonts/borders/
#    column layout match the rest of the workbook exactly) right after
#    "总计", then overwrite it with the 2022-Q4 fund holdings.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($null, $summary)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template has 26 data rows (rows 2-27); 2022-Q4 needs 29, so extend
# the sheet by 3 rows first, copying column A's formatting (bold/border/
# centered) down so the new rows match the rest of the index column.
for ($r = 28; $r -le 30; $r++) {
    $q4.Cells.Item(27, 1).Copy()
    $q4.Cells.Item($r, 1).PasteSpecial(-4122)
}

# Fund holdings for 2022-Q4 (index, code, name, fund size, stock
# position, position ratio, held market value, position rank).
$dataRows = @(
    @("0", "090018", "大成新锐产业混合", "88.75", "93.33", "9.37", "8.3159", "2"),
    @("1", "001300", "大成睿景灵活配置混合A", "32.95", "92.29", "9.19", "3.0281", "2"),
    @("2", "013435", "大成景气精选六个月持有混合A", "30.45", "91.16", "7.69", "2.3416", "3"),
    @("3", "001301", "大成睿景灵活配置混合C", "19.30", "92.29", "9.19", "1.7737", "2"),
    @("4", "002258", "大成国企改革灵活配置混合", "16.71", "93.37", "9.38", "1.5674", "2"),
    @("5", "014224", "大成聚优成长混合A", "16.55", "90.21", "7.85", "1.2992", "3"),
    @("6", "010826", "大成产业趋势混合A", "11.37", "93.99", "9.59", "1.0904", "2"),
    @("7", "012519", "大成核心趋势混合A", "9.51", "91.16", "9.45", "0.8987", "2"),
    @("8", "008934", "大成科技消费股票A", "10.09", "83.51", "4.55", "0.4591", "4"),
    @("9", "013436", "大成景气精选六个月持有混合C", "5.20", "91.16", "7.69", "0.3999", "3"),
    @("10", "012184", "大成创新趋势混合A", "7.64", "80.76", "4.65", "0.3553", "3"),
    @("11", "010827", "大成产业趋势混合C", "3.42", "93.99", "9.59", "0.3280", "2"),
    @("12", "014225", "大成聚优成长混合C", "3.26", "90.21", "7.85", "0.2559", "3"),
    @("13", "012520", "大成核心趋势混合C", "2.45", "91.16", "9.45", "0.2315", "2"),
    @("14", "008935", "大成科技消费股票C", "3.91", "83.51", "4.55", "0.1779", "4"),
    @("15", "014185", "招商专精特新股票A", "3.30", "87.72", "4.73", "0.1561", "7"),
    @("16", "008274", "大成行业先锋混合A", "2.81", "76.92", "4.39", "0.1234", "3"),
    @("17", "014186", "招商专精特新股票C", "2.50", "87.72", "4.73", "0.1182", "7"),
    @("18", "217013", "招商中小盘精选混合", "2.52", "86.61", "4.28", "0.1079", "7"),
    @("19", "002945", "大成盛世精选灵活配置混合", "1.24", "69.65", "4.66", "0.0578", "2"),
    @("20", "001531", "招商安益灵活配置混合", "0.55", "83.83", "3.92", "0.0216", "8"),
    @("21", "008275", "大成行业先锋混合C", "0.46", "76.92", "4.39", "0.0202", "3"),
    @("22", "015710", "华夏高端装备龙头混合A", "0.34", "94.44", "4.83", "0.0164", "7"),
    @("23", "620002", "金元顺安成长动力混合", "0.35", "72.02", "3.07", "0.0107", "10"),
    @("24", "012185", "大成创新趋势混合C", "0.15", "80.76", "4.65", "0.0070", "3"),
    @("25", "015711", "华夏高端装备龙头混合C", "0.13", "94.44", "4.83", "0.0063", "7"),
    @("26", "159620", "华夏中证智选500成长创新策略ETF", "0.32", "95.01", "1.51", "0.0048", "9"),
    @("27", "005966", "安信中证500指数增强C", "0.16", "88.79", "0.88", "0.0014", "7"),
    @("28", "005965", "安信中证500指数增强A", "0.12", "88.79", "0.88", "0.0011", "7")
)

for ($i = 0; $i -lt $dataRows.Length; $i++) {
    $row = $i + 2
    $vals = $dataRows[$i]

    $q4.Cells.Item($row, 1).Value = [int]$vals[0]

    # Columns B, D, E, F, G look numeric but must stay text (matches the
    # source data, and preserves the fund codes' leading zeros); the
    # leading apostrophe forces text entry the same way typing it in
    # Excel would.
    $q4.Cells.Item($row, 2).Value = "'" + $vals[1]
    $q4.Cells.Item($row, 3).Value = $vals[2]
    $q4.Cells.Item($row, 4).Value = "'" + $vals[3]
    $q4.Cells.Item($row, 5).Value = "'" + $vals[4]
    $q4.Cells.Item($row, 6).Value = "'" + $vals[5]
    $q4.Cells.Item($row, 7).Value = "'" + $vals[6]
    $q4.Cells.Item($row, 8).Value = [int]$vals[7]
}

Write-Output "2022-Q4 sheet added with $($dataRows.Length) holdings."
